$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so strings like "303.49" are
# not auto-converted to numeric values by Excel's smart cell parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.424.23'
$ws.Range('E2').Value = '  -1.40%  '

$ws.Range('D3').Value = '2.286.02'
$ws.Range('E3').Value = '  -0.77%  '

$ws.Range('D5').Value = '303.49'
$ws.Range('E5').Value = '  +0.98%  '

$ws.Range('D6').Value = '95.09'
$ws.Range('E6').Value = '  -3.19%  '

$ws.Range('E7').Value = '  -2.95%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('E9').Value = '  -3.50%  '

$ws.Range('D10').Value = '34.91'
$ws.Range('E10').Value = '  -3.95%  '

$ws.Range('D11').Value = '0.0778'
$ws.Range('E11').Value = '  -1.42%  '

$ws.Range('E12').Value = '  +1.60%  '

$ws.Range('D13').Value = '17.92'
$ws.Range('E13').Value = '  -0.26%  '

$ws.Range('E14').Value = '  -0.38%  '

$ws.Range('D15').Value = '2.638.51'
$ws.Range('E15').Value = '  -0.91%  '

$ws.Range('D16').Value = '2.280.29'
$ws.Range('E16').Value = '  -1.56%  '

$ws.Range('E17').Value = '  -1.77%  '

$ws.Range('D18').Value = '42.349.23'
$ws.Range('E18').Value = '  -1.47%  '

$ws.Range('D19').Value = '12.67'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('D20').Value = '0.0₃0886'
$ws.Range('E20').Value = '  -2.65%  '

$ws.Range('D21').Value = '5.95'
$ws.Range('E21').Value = '  -2.82%  '

$ws.Range('E22').Value = '  -1.99%  '

$ws.Range('D23').Value = '235.77'
$ws.Range('E23').Value = '  -2.58%  '

$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  +0.13%  '

$ws.Range('E25').Value = '  +0.04%  '

$ws.Range('D26').Value = '2.39'
$ws.Range('E26').Value = '  -2.00%  '

$ws.Range('E27').Value = '  -2.31%  '

$ws.Range('E28').Value = '  +16.99%  '

$ws.Range('D29').Value = '167.56'
$ws.Range('E29').Value = '  +0.44%  '

$ws.Range('D30').Value = '8.94'
$ws.Range('E30').Value = '  -1.82%  '

$ws.Range('D31').Value = '32.28'
$ws.Range('E31').Value = '  -3.08%  '

$ws.Range('E32').Value = '  +0.03%  '

$ws.Range('D33').Value = '17.62'
$ws.Range('E33').Value = '  -0.63%  '

$ws.Range('D34').Value = '4.92'
$ws.Range('E34').Value = '  -2.18%  '

$ws.Range('D35').Value = '4.46'
$ws.Range('E35').Value = '  -6.74%  '

$ws.Range('E36').Value = '  -2.46%  '

$ws.Range('E37').Value = '  -1.03%  '

$ws.Range('E38').Value = '  -1.07%  '

$ws.Range('E39').Value = '  -2.71%  '

$ws.Range('E40').Value = '  -2.54%  '

$ws.Range('D41').Value = '2.65'
$ws.Range('E41').Value = '  -4.72%  '

$ws.Range('D42').Value = '1.984.22'
$ws.Range('E42').Value = '  -0.66%  '

$ws.Range('D43').Value = '0.0275'
$ws.Range('E43').Value = '  -4.05%  '

$ws.Range('D44').Value = '10.02'
$ws.Range('E44').Value = '  -1.99%  '

$ws.Range('D45').Value = '17.79'
$ws.Range('E45').Value = '  +1.52%  '

$ws.Range('D47').Value = '2.74'
$ws.Range('E47').Value = '  -2.38%  '

$ws.Range('D48').Value = '2.96'
$ws.Range('E48').Value = '  +6.76%  '

$ws.Range('D49').Value = '53.16'
$ws.Range('E49').Value = '  -0.88%  '

$ws.Range('D50').Value = '2.504.57'
$ws.Range('E50').Value = '  -1.02%  '

$ws.Range('D51').Value = '70.31'
$ws.Range('E51').Value = '  -3.62%  '

# Restore the original cell style (no explicit style / numFmt) now that the
# text values are safely stored, so formatting matches the source workbook.
$ws.Range("D2:D51").Style = "Normal"
